$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the employee name (shared string "test" -> "曾閔歆")
$ws.Range("B1").Value = "曾閔歆"

# Update existing data rows 3 and 4
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = 18
$ws.Range("C3").Value = 16
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 21
$ws.Range("F3").Value = 20

$ws.Range("A4").Value = 11
$ws.Range("B4").Value = 25
$ws.Range("C4").Value = 16
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 21
$ws.Range("F4").Value = 25

# New row 5
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = 28
$ws.Range("C5").Value = 16
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 21
$ws.Range("F5").Value = 30

# New row 6 (last data row)
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = 29
$ws.Range("C6").Value = 15
$ws.Range("D6").Value = 30
$ws.Range("E6").Value = 19
$ws.Range("F6").Value = 3

# Trailer rows 7-15: blank placeholder cells in columns C and E only.
# Touching Borders materializes the (valueless) cell in the sheet XML,
# matching the style used for the other borderless C/E cells above.
for ($r = 7; $r -le 15; $r++) {
    $ws.Cells.Item($r, 3).Borders.LineStyle = -4142
    $ws.Cells.Item($r, 5).Borders.LineStyle = -4142
}

# Final selection, matching the saved workbook view
[void]$ws.Range("E6").Select()
